$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.089.78"
$ws.Range("E2").Value = "'  +1.71%  "
$ws.Range("D3").Value = "'1.856.34"
$ws.Range("E3").Value = "'  +3.06%  "
$ws.Range("E4").Value = "'  +0.31%  "
$ws.Range("D5").Value = "'237.21"
$ws.Range("E5").Value = "'  +3.51%  "
$ws.Range("D6").Value = "'0.623"
$ws.Range("E6").Value = "'  +2.09%  "
$ws.Range("E7").Value = "'  +0.16%  "
$ws.Range("D8").Value = "'42.28"
$ws.Range("E9").Value = "'  +3.02%  "
$ws.Range("E10").Value = "'  +2.55%  "
$ws.Range("D11").Value = "'0.0992"
$ws.Range("E11").Value = "'  +0.40%  "
$ws.Range("D12").Value = "'2.126.14"
$ws.Range("E12").Value = "'  +3.13%  "
$ws.Range("B13").Value = "'Chainlink"
$ws.Range("C13").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.44"
$ws.Range("E13").Value = "'  +3.11%  "
$ws.Range("B14").Value = "'WrappedEther"
$ws.Range("C14").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.860.15"
$ws.Range("E14").Value = "'  +3.44%  "
$ws.Range("D15").Value = "'0.678"
$ws.Range("E15").Value = "'  +2.89%  "
$ws.Range("E16").Value = "'  +3.17%  "
$ws.Range("D17").Value = "'35.041.32"
$ws.Range("E17").Value = "'  +1.98%  "
$ws.Range("D18").Value = "'70.36"
$ws.Range("D19").Value = "'0.0₃0795"
$ws.Range("E19").Value = "'  +2.08%  "
$ws.Range("D20").Value = "'240.68"
$ws.Range("E20").Value = "'  +0.58%  "
$ws.Range("D21").Value = "'12.13"
$ws.Range("E21").Value = "'  +3.03%  "
$ws.Range("D22").Value = "'4.75"
$ws.Range("E22").Value = "'  +1.60%  "
$ws.Range("E24").Value = "'  +1.85%  "
$ws.Range("D25").Value = "'171.24"
$ws.Range("E25").Value = "'  -0.90%  "
$ws.Range("E26").Value = "'  +27.29%  "
$ws.Range("E27").Value = "'  +2.84%  "
$ws.Range("D28").Value = "'17.68"
$ws.Range("E28").Value = "'  +3.02%  "
$ws.Range("E29").Value = "'  +2.32%  "
$ws.Range("B30").Value = "'BinanceUSD"
$ws.Range("C30").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D30").Value = "'1.01"
$ws.Range("E30").Value = "'  +0.29%  "
$ws.Range("B31").Value = "'Hedera"
$ws.Range("C31").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0558"
$ws.Range("E31").Value = "'  +3.06%  "
$ws.Range("D32").Value = "'4.00"
$ws.Range("E32").Value = "'  +0.45%  "
$ws.Range("D33").Value = "'4.01"
$ws.Range("E33").Value = "'  +2.98%  "
$ws.Range("E34").Value = "'  +13.15%  "
$ws.Range("E35").Value = "'  +22.97%  "
$ws.Range("E36").Value = "'  +5.63%  "
$ws.Range("E37").Value = "'  +13.23%  "
$ws.Range("D38").Value = "'1.08"
$ws.Range("E38").Value = "'  +13.07%  "
$ws.Range("D39").Value = "'91.80"
$ws.Range("E39").Value = "'  +1.62%  "
$ws.Range("E40").Value = "'  +7.27%  "
$ws.Range("D41").Value = "'1.353.49"
$ws.Range("E41").Value = "'  +2.30%  "
$ws.Range("D42").Value = "'14.85"
$ws.Range("E42").Value = "'  +4.50%  "
$ws.Range("E43").Value = "'  +6.26%  "
$ws.Range("D44").Value = "'12.61"
$ws.Range("E44").Value = "'  +55.18%  "
$ws.Range("D45").Value = "'2.41"
$ws.Range("E45").Value = "'  +1.02%  "
$ws.Range("D46").Value = "'2.74"
$ws.Range("E46").Value = "'  +1.45%  "
$ws.Range("E47").Value = "'  +7.54%  "
$ws.Range("E48").Value = "'  +5.42%  "
$ws.Range("D49").Value = "'2.037.95"
$ws.Range("E49").Value = "'  +2.65%  "
$ws.Range("E50").Value = "'  +3.40%  "
$ws.Range("D51").Value = "'3.42"
$ws.Range("E51").Value = "'  +17.87%  "
